$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated "K" column (G) values for rows 2-12, replacing old Strike# values
$values = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 2
    7  = 6
    8  = 3
    9  = 3
    10 = 3
    11 = 3
    12 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
